$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '65.026.58'
$ws.Range("E2").Value = '  +1.40%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.138.68'
$ws.Range("E3").Value = '  +3.09%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.21%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '566.56'
$ws.Range("E5").Value = '  +1.53%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '148.65'
$ws.Range("E6").Value = '  +6.10%  '

$ws.Range("E7").Value = '  -0.03%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.129.39'
$ws.Range("E8").Value = '  +3.02%  '

$ws.Range("E9").Value = '  +3.25%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.06'
$ws.Range("E10").Value = '  +15.91%  '

$ws.Range("E11").Value = '  +2.65%  '

$ws.Range("E12").Value = '  +1.41%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '36.03'
$ws.Range("E13").Value = '  +3.27%  '

$ws.Range("E14").Value = '  +1.25%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.644.41'
$ws.Range("E15").Value = '  +2.76%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.019.75'
$ws.Range("E16").Value = '  +1.19%  '

$ws.Range("B17").Value = 'BitcoinCash'
$ws.Range("C17").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '534.70'
$ws.Range("E17").Value = '  +11.28%  '

$ws.Range("B18").Value = 'TRON'
$ws.Range("C18").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.113'
$ws.Range("E18").Value = '  +2.43%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.138.43'
$ws.Range("E19").Value = '  +2.81%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.77'
$ws.Range("E20").Value = '  +3.33%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.90'

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.706'
$ws.Range("E22").Value = '  +5.26%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.44'
$ws.Range("E23").Value = '  +4.70%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.81'
$ws.Range("E24").Value = '  +3.22%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '78.84'
$ws.Range("E25").Value = '  +1.71%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  +0.02%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.76'
$ws.Range("E27").Value = '  +14.77%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.81'
$ws.Range("E28").Value = '  +2.91%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.14'
$ws.Range("E29").Value = '  +3.85%  '

$ws.Range("E30").Value = '  +0.10%  '

$ws.Range("B31").Value = 'EthereumClassic'
$ws.Range("C31").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '26.24'
$ws.Range("E31").Value = '  +0.86%  '

$ws.Range("B32").Value = 'Stacks'
$ws.Range("C32").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.64'
$ws.Range("E32").Value = '  +3.42%  '

$ws.Range("E33").Value = '  +4.66%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '552.34'
$ws.Range("E34").Value = '  +10.71%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.41'
$ws.Range("E35").Value = '  +1.74%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.08'
$ws.Range("E36").Value = '  +4.61%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0445'
$ws.Range("E37").Value = '  +9.42%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '52.99'
$ws.Range("E38").Value = '  +0.88%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0819'
$ws.Range("E39").Value = '  +4.59%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.87'
$ws.Range("E40").Value = '  +11.98%  '

$ws.Range("B41").Value = 'Kaspa'
$ws.Range("C41").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.122'
$ws.Range("E41").Value = '  +2.81%  '

$ws.Range("B42").Value = 'Maker'
$ws.Range("C42").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.065.61'
$ws.Range("E42").Value = '  +7.32%  '

$ws.Range("E43").Value = '  +0.63%  '

$ws.Range("E44").Value = '  +7.55%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.18'
$ws.Range("E45").Value = '  +7.55%  '

$ws.Range("E46").Value = '  +0.12%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '25.05'
$ws.Range("E47").Value = '  +2.76%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '120.04'
$ws.Range("E48").Value = '  +2.13%  '

$ws.Range("E49").Value = '  -1.38%  '

$ws.Range("E50").Value = '  +2.69%  '

$ws.Range("E51").Value = '  +3.53%  '
